# Insert a new data row at row 125 (pushing the existing rows 125-184 down
# to 126-185) and populate it with the new weekly record. This mirrors the
# OOXML diff: <dimension> grows from A1:R184 to A1:R185, and every row from
# 125 onward in the "before" file reappears one row lower in the "after"
# file, with a brand-new row 125 inserted in front of them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 125..184 down to 126..185, creating a blank row 125.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the new record's values.
$ws.Cells.Item(125, 1).Value = 8
$ws.Cells.Item(125, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(125, 3).Value = "Coquimbo"
$ws.Cells.Item(125, 4).Value = 44603
$ws.Cells.Item(125, 5).Value = 4
$ws.Cells.Item(125, 6).Value = 100112031
$ws.Cells.Item(125, 7).Value = "Poroto verde"
$ws.Cells.Item(125, 8).Value = "Magnum"
$ws.Cells.Item(125, 9).Value = "Primera"
$ws.Cells.Item(125, 10).Value = 480
$ws.Cells.Item(125, 11).Value = 36000
$ws.Cells.Item(125, 12).Value = 37000
$ws.Cells.Item(125, 13).Value = 36500
$ws.Cells.Item(125, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(125, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(125, 16).Value = 1460
$ws.Cells.Item(125, 17).Value = 25
$ws.Cells.Item(125, 18).Value = "Hortaliza"
